$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9473599791526794
$ws.Range("B1").Value = 1.99790096282959
$ws.Range("C1").Value = 2.950685024261475
$ws.Range("D1").Value = 3.582872152328491
$ws.Range("E1").Value = 1.974338173866272
